# Update to schema: split the single "methodTypes" column into two columns,
# "supportingMethods" followed by a new "supportingMethodTypes" column, on
# every worksheet that has a methodTypes column.
#
# For each sheet: rename the existing methodTypes cell to supportingMethods,
# then insert a new blank column immediately to its right and label it
# supportingMethodTypes (which pushes every later column one slot to the
# right, exactly as in the target diff).

$wb = $excel.ActiveWorkbook

$sheetsAndColumns = @{
    "InformationEntity" = 5   # E1: methodTypes
    "Method"             = 7   # G1: methodTypes
    "Document"           = 11  # K1: methodTypes
    "DataItem"           = 11  # K1: methodTypes
    "DataSet"            = 9   # I1: methodTypes
    "Statement"          = 18  # R1: methodTypes
    "StudyResult"        = 11  # K1: methodTypes
    "EvidenceLine"       = 12  # L1: methodTypes
}

foreach ($sheetName in $sheetsAndColumns.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $col = $sheetsAndColumns[$sheetName]

    # Rename the existing "methodTypes" header to "supportingMethods".
    $ws.Cells.Item(1, $col).Value = "supportingMethods"

    # Insert a new column right after it and label it "supportingMethodTypes".
    $ws.Columns.Item($col + 1).Insert()
    $ws.Cells.Item(1, $col + 1).Value = "supportingMethodTypes"
}
